# Automatic update of files.
# The edit re-shuffles the per-record data (id, taxon info, coordinates, ...)
# held in columns A-N/Q/R across rows 10-22 and 24-28 of the "Artfynd" sheet,
# while columns O/P/S and T..AY (shared metadata) stay in place.
#
# Mapping: destination row -> source row (values are copied from the
# ORIGINAL/"before" content of the source row into the destination row).
#   10 <- 12      12 <- 17      13 <- 18      14 <- 21      16 <- 19
#   17 <- 20      18 <- 13      19 <- 22      20 <- 14      21 <- 16
#   22 <- 10      24 <- 28      25 <- 27      27 <- 24      28 <- 25

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns that carry the per-record data which gets reshuffled.
$cols = @("A", "B", "D", "E", "F", "G", "H", "K", "L", "M", "N", "Q", "R")

$srcRows = @(10, 12, 13, 14, 16, 17, 18, 19, 20, 21, 22, 24, 25, 27, 28)

# Snapshot the original ("before") values of every involved row/column so
# that writing the destinations doesn't clobber values still needed as a
# source later on (the remapping contains cycles).
$snapshot = @{}
foreach ($r in $srcRows) {
    foreach ($c in $cols) {
        $addr = "$c$r"
        $snapshot[$addr] = $ws.Range($addr).Value()
    }
}

# destination row -> source row
$mapping = @{
    10 = 12
    12 = 17
    13 = 18
    14 = 21
    16 = 19
    17 = 20
    18 = 13
    19 = 22
    20 = 14
    21 = 16
    22 = 10
    24 = 28
    25 = 27
    27 = 24
    28 = 25
}

foreach ($dst in $mapping.Keys) {
    $src = $mapping[$dst]
    foreach ($c in $cols) {
        $val = $snapshot["$c$src"]
        $addr = "$c$dst"
        if ($null -eq $val) {
            $ws.Range($addr).Value = ""
        } else {
            $ws.Range($addr).Value = $val
        }
    }
}
